# Applies the commit "changed P inputs & corrected Med output":
#  - Column C (k) updated from 1 to 0.8 for every strut row (P input change).
#  - Column D (s.m) reduced from 4 to 3 for the "S2" rows (P input change).
#  - Strut names (column A) lose their leading "A-" prefix; rows 8-11 are
#    renumbered from "L3" to "L4"/"L5" to reflect the corrected level layout.
#  - All of the downstream "Med"/capacity-check outputs (columns X, Y, Z, AA,
#    AB..BA) are recalculated to match the corrected member design outputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("A2", "L1-S1"),
    @("C2", 0.8),
    @("AD2", 3876),
    @("AE2", 1.12),
    @("AG2", 0.58),
    @("AH2", 2839),
    @("AL2", 71505),
    @("AM2", 0.26),
    @("AO2", 0.99),
    @("AP2", 4832),
    @("AS2", 33032),
    @("AT2", 0.27),
    @("AV2", 0.98),
    @("AW2", 2410),
    @("AY2", 71505),
    @("AZ2", 41),
    @("A3", "L1-S2"),
    @("C3", 0.8),
    @("D3", 3),
    @("X3", 219),
    @("Y3", 2242),
    @("Z3", 1121),
    @("AA3", "305 x 127 x 37"),
    @("AC3", 3351),
    @("AD3", 2867),
    @("AE3", 1.08),
    @("AG3", 0.61),
    @("AH3", 2041),
    @("AJ3", 3351),
    @("AK3", 236000),
    @("AL3", 94355),
    @("AM3", 0.19),
    @("AO3", 1),
    @("AP3", 3360),
    @("AR3", 1676),
    @("AS3", 10881),
    @("AT3", 0.39),
    @("AV3", 0.95),
    @("AW3", 1600),
    @("AY3", 94355),
    @("AZ3", 20),
    @("BA3", 563),
    @("A4", "L2-S1"),
    @("C4", 0.8),
    @("AD4", 28188),
    @("AE4", 0.66),
    @("AG4", 0.87),
    @("AH4", 10632),
    @("AL4", 184463),
    @("AM4", 0.26),
    @("AO4", 0.99),
    @("AP4", 12124),
    @("AS4", 167752),
    @("AT4", 0.19),
    @("AV4", 1),
    @("AW4", 6153),
    @("AY4", 184463),
    @("AZ4", 125),
    @("BA4", 2440),
    @("A5", "L2-S2"),
    @("C5", 0.8),
    @("D5", 3),
    @("Y5", 6866),
    @("Z5", 3433),
    @("AA5", "533 x 210 x 92"),
    @("AB5", 355),
    @("AC5", 8307),
    @("AD5", 22070),
    @("AE5", 0.61),
    @("AG5", 0.88),
    @("AH5", 7350),
    @("AI5", 355),
    @("AJ5", 8307),
    @("AK5", 585000),
    @("AL5", 233889),
    @("AM5", 0.19),
    @("AO5", 1),
    @("AP5", 8328),
    @("AQ5", 355),
    @("AR5", 4154),
    @("AS5", 77399),
    @("AT5", 0.23),
    @("AW5", 4125),
    @("AY5", 233889),
    @("AZ5", 63),
    @("BA5", 1723),
    @("A6", "L3-S1"),
    @("C6", 0.8),
    @("AD6", 77930),
    @("AE6", 0.5),
    @("AG6", 0.92),
    @("AH6", 18220),
    @("AL6", 296384),
    @("AM6", 0.26),
    @("AO6", 0.99),
    @("AP6", 19480),
    @("AS6", 362708),
    @("AT6", 0.16),
    @("AV6", 1.01),
    @("AW6", 9942),
    @("AY6", 296384),
    @("AZ6", 210),
    @("A7", "L3-S2"),
    @("C7", 0.8),
    @("D7", 3),
    @("Y7", 11574),
    @("Z7", 5787),
    @("AA7", "686 x 254 x 152"),
    @("AC7", 13386),
    @("AD7", 59972),
    @("AE7", 0.47),
    @("AG7", 0.93),
    @("AH7", 12484),
    @("AJ7", 13386),
    @("AK7", 970000),
    @("AL7", 387816),
    @("AM7", 0.19),
    @("AO7", 1),
    @("AP7", 13428),
    @("AR7", 6693),
    @("AS7", 187183),
    @("AT7", 0.19),
    @("AW7", 6709),
    @("AY7", 387816),
    @("AZ7", 106),
    @("BA7", 2904),
    @("A8", "L4-S1"),
    @("C8", 0.8),
    @("AD8", 49743),
    @("AE8", 0.59),
    @("AG8", 0.89),
    @("AH8", 15478),
    @("AL8", 260113),
    @("AM8", 0.26),
    @("AO8", 0.99),
    @("AP8", 17096),
    @("AS8", 264583),
    @("AT8", 0.18),
    @("AV8", 1),
    @("AW8", 8696),
    @("AY8", 260113),
    @("AZ8", 185),
    @("BA8", 3614),
    @("A9", "L4-S2"),
    @("C9", 0.8),
    @("D9", 3),
    @("Y9", 10170),
    @("Z9", 5085),
    @("AA9", "686 x 254 x 140"),
    @("AC9", 12282),
    @("AD9", 54374),
    @("AE9", 0.48),
    @("AG9", 0.93),
    @("AH9", 11444),
    @("AJ9", 12282),
    @("AK9", 890000),
    @("AL9", 355831),
    @("AM9", 0.19),
    @("AO9", 1),
    @("AP9", 12320),
    @("AR9", 6141),
    @("AS9", 167752),
    @("AT9", 0.19),
    @("AW9", 6153),
    @("AY9", 355831),
    @("AZ9", 93),
    @("BA9", 2552),
    @("A10", "L5-S1"),
    @("C10", 0.8),
    @("AD10", 24457),
    @("AE10", 0.67),
    @("AG10", 0.86),
    @("AH10", 9451),
    @("AL10", 164773),
    @("AM10", 0.26),
    @("AO10", 0.99),
    @("AP10", 10830),
    @("AS10", 141845),
    @("AT10", 0.2),
    @("AV10", 1),
    @("AW10", 5490),
    @("AY10", 164773),
    @("AZ10", 103),
    @("A11", "L5-S2"),
    @("C11", 0.8),
    @("D11", 3),
    @("Y11", 5664),
    @("Z11", 2832),
    @("AA11", "406 x 178 x 74"),
    @("AB11", 355),
    @("AC11", 6710),
    @("AD11", 10915),
    @("AE11", 0.78),
    @("AG11", 0.8),
    @("AH11", 5399),
    @("AI11", 355),
    @("AJ11", 6710),
    @("AK11", 472500),
    @("AL11", 188910),
    @("AM11", 0.19),
    @("AO11", 1),
    @("AP11", 6726),
    @("AQ11", 355),
    @("AR11", 3355),
    @("AS11", 50196),
    @("AT11", 0.26),
    @("AV11", 0.99),
    @("AW11", 3311),
    @("AY11", 188910),
    @("AZ11", 52),
    @("BA11", 1421)
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $val = $pair[1]
    $ws.Range($addr).Value = $val
}

Write-Host ("Applied {0} cell updates" -f $updates.Count)
